$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 309
$ws.Range("I5").Value = 257
$ws.Range("K5").Value = 257
$ws.Range("M5").Value = -142
$ws.Range("H33").Value = 286.30768
$ws.Range("I33").Value = 130.5
$ws.Range("K33").Value = 130.5
$ws.Range("M33").Value = 98.5
$ws.Range("H137").Value = 1088.8636
$ws.Range("I137").Value = 856.2059
$ws.Range("J137").Value = 1336.0625
$ws.Range("K137").Value = 2568.6177
$ws.Range("L137").Value = 4008.1875
$ws.Range("M137").Value = -18.61770000000024
$ws.Range("N137").Value = -9108.1875
$ws.Range("H138").Value = 1516.35
$ws.Range("I138").Value = 1019.7143
$ws.Range("J138").Value = 1553.7312
$ws.Range("K138").Value = 3059.1429
$ws.Range("L138").Value = 4661.1936
$ws.Range("M138").Value = 2080.8571
$ws.Range("N138").Value = -14941.1936

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2952.015
$ws.Range("I32").Value = 2696.6924
$ws.Range("K32").Value = 2696.6924
$ws.Range("M32").Value = -2409.6924
$ws.Range("H61").Value = 1276.409
$ws.Range("I61").Value = 1174.8286
$ws.Range("J61").Value = 1671.4445
$ws.Range("K61").Value = 1174.8286
$ws.Range("L61").Value = 1671.4445
$ws.Range("M61").Value = -962.8286000000001
$ws.Range("N61").Value = -2095.4445
$ws.Range("H74").Value = 1045.4
$ws.Range("I74").Value = 570
$ws.Range("K74").Value = 570
$ws.Range("M74").Value = 304
$ws.Range("H77").Value = 1045.4
$ws.Range("I77").Value = 570
$ws.Range("K77").Value = 2850
$ws.Range("M77").Value = 1518
$ws.Range("H93").Value = 32900
$ws.Range("J93").Value = 32900
$ws.Range("L93").Value = 32900
$ws.Range("N93").Value = -37892
$ws.Range("H132").Value = 1887.3889
$ws.Range("I132").Value = 1604.24
$ws.Range("K132").Value = 4812.72
$ws.Range("M132").Value = -2282.72
$ws.Range("H136").Value = 1276.409
$ws.Range("I136").Value = 1174.8286
$ws.Range("J136").Value = 1671.4445
$ws.Range("K136").Value = 3524.4858
$ws.Range("L136").Value = 5014.333500000001
$ws.Range("M136").Value = -974.4858000000004
$ws.Range("N136").Value = -10114.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 491
$ws.Range("I5").Value = 477
$ws.Range("J5").Value = 505
$ws.Range("K5").Value = 477
$ws.Range("L5").Value = 505
$ws.Range("M5").Value = -364
$ws.Range("N5").Value = -731
$ws.Range("H134").Value = 3832.5
$ws.Range("I134").Value = 954.2059
$ws.Range("J134").Value = 13618.7
$ws.Range("K134").Value = 2862.6177
$ws.Range("L134").Value = 40856.10000000001
$ws.Range("M134").Value = -327.6177000000002
$ws.Range("N134").Value = -45926.10000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1890.8667
$ws.Range("I31").Value = 2014.2609
$ws.Range("J31").Value = 1485.4286
$ws.Range("K31").Value = 2014.2609
$ws.Range("L31").Value = 1485.4286
$ws.Range("M31").Value = -1719.2609
$ws.Range("N31").Value = -2075.4286
$ws.Range("H34").Value = 1890.8667
$ws.Range("I34").Value = 2014.2609
$ws.Range("J34").Value = 1485.4286
$ws.Range("K34").Value = 2014.2609
$ws.Range("L34").Value = 1485.4286
$ws.Range("M34").Value = -1812.2609
$ws.Range("N34").Value = -1889.4286
$ws.Range("H41").Value = 13350
$ws.Range("I41").Value = 2733.3333
$ws.Range("J41").Value = 23966.666
$ws.Range("K41").Value = 2733.3333
$ws.Range("L41").Value = 23966.666
$ws.Range("M41").Value = -2305.3333
$ws.Range("N41").Value = -24822.666
$ws.Range("H43").Value = 5637.25
$ws.Range("J43").Value = 5637.25
$ws.Range("L43").Value = 5637.25
$ws.Range("N43").Value = -6005.25
$ws.Range("H50").Value = 26886
$ws.Range("J50").Value = 26886
$ws.Range("L50").Value = 26886
$ws.Range("N50").Value = -28136
$ws.Range("H51").Value = 18800
$ws.Range("J51").Value = 22250
$ws.Range("L51").Value = 22250
$ws.Range("N51").Value = -23722
$ws.Range("H58").Value = 796.67566
$ws.Range("I58").Value = 790.087
$ws.Range("J58").Value = 807.5
$ws.Range("K58").Value = 790.087
$ws.Range("L58").Value = 807.5
$ws.Range("M58").Value = -587.087
$ws.Range("N58").Value = -1213.5
$ws.Range("H59").Value = 25000
$ws.Range("J59").Value = 25000
$ws.Range("L59").Value = 25000
$ws.Range("N59").Value = -27290
$ws.Range("H60").Value = 3142.8572
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H61").Value = 18800
$ws.Range("J61").Value = 22250
$ws.Range("L61").Value = 22250
$ws.Range("N61").Value = -22946
$ws.Range("H99").Value = 2170.2
$ws.Range("I99").Value = 2200.2222
$ws.Range("J99").Value = 1900
$ws.Range("K99").Value = 2200.2222
$ws.Range("L99").Value = 1900
$ws.Range("M99").Value = -702.2222000000002
$ws.Range("N99").Value = -4896
$ws.Range("H101").Value = 5637.25
$ws.Range("J101").Value = 5637.25
$ws.Range("L101").Value = 5637.25
$ws.Range("N101").Value = -12127.25
$ws.Range("H122").Value = 1894.7142
$ws.Range("I122").Value = 1353
$ws.Range("J122").Value = 2111.4
$ws.Range("K122").Value = 4059
$ws.Range("L122").Value = 6334.200000000001
$ws.Range("M122").Value = -1609
$ws.Range("N122").Value = -11234.2
$ws.Range("H126").Value = 2170.2
$ws.Range("I126").Value = 2200.2222
$ws.Range("J126").Value = 1900
$ws.Range("K126").Value = 6600.6666
$ws.Range("L126").Value = 5700
$ws.Range("M126").Value = -4130.6666
$ws.Range("N126").Value = -10640
$ws.Range("H134").Value = 2121.4
$ws.Range("I134").Value = 2416.2856
$ws.Range("J134").Value = 1433.3334
$ws.Range("K134").Value = 7248.8568
$ws.Range("L134").Value = 4300.0002
$ws.Range("M134").Value = -4713.8568
$ws.Range("N134").Value = -9370.0002
$ws.Range("H136").Value = 796.67566
$ws.Range("I136").Value = 790.087
$ws.Range("J136").Value = 807.5
$ws.Range("K136").Value = 2370.261
$ws.Range("L136").Value = 2422.5
$ws.Range("M136").Value = 179.739
$ws.Range("N136").Value = -7522.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 4503.2
$ws.Range("J42").Value = 4503.2
$ws.Range("L42").Value = 13509.6
$ws.Range("N42").Value = -14577.6
$ws.Range("H92").Value = 281.9524
$ws.Range("I92").Value = 387.6
$ws.Range("J92").Value = 248.9375
$ws.Range("K92").Value = 1162.8
$ws.Range("L92").Value = 746.8125
$ws.Range("M92").Value = 85.19999999999982
$ws.Range("N92").Value = -3242.8125
$ws.Range("H113").Value = 714.069
$ws.Range("I113").Value = 700
$ws.Range("J113").Value = 714.5714
$ws.Range("K113").Value = 2100
$ws.Range("L113").Value = 2143.7142
$ws.Range("M113").Value = 70
$ws.Range("N113").Value = -6483.7142
$ws.Range("H131").Value = 22728572
$ws.Range("J131").Value = 1603.7059
$ws.Range("L131").Value = 4811.1177
$ws.Range("N131").Value = -14891.1177
$ws.Range("H139").Value = 1829.7693
$ws.Range("I139").Value = 1954.619
$ws.Range("K139").Value = 5863.857
$ws.Range("M139").Value = -723.857

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 196
$ws.Range("I2").Value = 261
$ws.Range("K2").Value = 261
$ws.Range("M2").Value = -148
$ws.Range("H113").Value = 1608.5
$ws.Range("I113").Value = 1555.75
$ws.Range("J113").Value = 1925
$ws.Range("K113").Value = 1555.75
$ws.Range("L113").Value = 1925
$ws.Range("M113").Value = 614.25
$ws.Range("N113").Value = -6265
$ws.Range("H132").Value = 2246.5789
$ws.Range("I132").Value = 1917.0625
$ws.Range("J132").Value = 4004
$ws.Range("K132").Value = 5751.1875
$ws.Range("L132").Value = 12012
$ws.Range("M132").Value = -3221.1875
$ws.Range("N132").Value = -17072

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1569.625
$ws.Range("I7").Value = 1308.6666
$ws.Range("K7").Value = 1308.6666
$ws.Range("M7").Value = -1196.6666
$ws.Range("H82").Value = 2065.2666
$ws.Range("I82").Value = 2075
$ws.Range("J82").Value = 2054.1428
$ws.Range("K82").Value = 2075
$ws.Range("L82").Value = 2054.1428
$ws.Range("M82").Value = -1714
$ws.Range("N82").Value = -2776.1428
$ws.Range("H85").Value = 2065.2666
$ws.Range("I85").Value = 2075
$ws.Range("J85").Value = 2054.1428
$ws.Range("K85").Value = 2075
$ws.Range("L85").Value = 2054.1428
$ws.Range("M85").Value = -827
$ws.Range("N85").Value = -4550.1428
$ws.Range("H126").Value = 1569.625
$ws.Range("I126").Value = 1308.6666
$ws.Range("K126").Value = 3925.9998
$ws.Range("M126").Value = -1455.9998
$ws.Range("H132").Value = 18117.393
$ws.Range("I132").Value = 1259.8055
$ws.Range("K132").Value = 3779.4165
$ws.Range("M132").Value = -1249.4165
$ws.Range("H135").Value = 33889.668
$ws.Range("J135").Value = 33889.668
$ws.Range("L135").Value = 33889.668
$ws.Range("N135").Value = -44029.668
$ws.Range("H136").Value = 1338.4
$ws.Range("I136").Value = 1033.1428
$ws.Range("J136").Value = 2050.6667
$ws.Range("K136").Value = 3099.4284
$ws.Range("L136").Value = 6152.000100000001
$ws.Range("M136").Value = -549.4284000000002
$ws.Range("N136").Value = -11252.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 467.70587
$ws.Range("I107").Value = 436.73334
$ws.Range("J107").Value = 700
$ws.Range("K107").Value = 1310.20002
$ws.Range("L107").Value = 2100
$ws.Range("M107").Value = 609.79998
$ws.Range("N107").Value = -5940
$ws.Range("H126").Value = 58824988
$ws.Range("I126").Value = 71429810
$ws.Range("K126").Value = 214289430
$ws.Range("M126").Value = -214286960
$ws.Range("H132").Value = 2677.0789
$ws.Range("I132").Value = 2478.3142
$ws.Range("J132").Value = 4996
$ws.Range("K132").Value = 7434.942599999999
$ws.Range("L132").Value = 14988
$ws.Range("M132").Value = -4904.942599999999
$ws.Range("N132").Value = -20048
$ws.Range("H136").Value = 693.0769
$ws.Range("I136").Value = 356.75
$ws.Range("J136").Value = 1231.2
$ws.Range("K136").Value = 1070.25
$ws.Range("L136").Value = 3693.6
$ws.Range("M136").Value = 1479.75
$ws.Range("N136").Value = -8793.6
